$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A6").Value = "Jeffery"
$ws.Range("B6").Value = "Chen"
